$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.403.44"
Set-TextValue $ws.Range("E2") "  -1.52%  "
Set-TextValue $ws.Range("D3") "1.711.25"
Set-TextValue $ws.Range("E3") "  -1.66%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "224.59"
Set-TextValue $ws.Range("E5") "  -1.38%  "
Set-TextValue $ws.Range("D6") "0.5334"
Set-TextValue $ws.Range("E6") "  -2.37%  "
Set-TextValue $ws.Range("D7") "1.004"
Set-TextValue $ws.Range("E7") "  +0.01%  "
Set-TextValue $ws.Range("D8") "0.2670"
Set-TextValue $ws.Range("E8") "  -3.93%  "
Set-TextValue $ws.Range("D9") "0.06610"
Set-TextValue $ws.Range("E9") "  -2.15%  "
Set-TextValue $ws.Range("D10") "20.92"
Set-TextValue $ws.Range("E10") "  -4.94%  "
Set-TextValue $ws.Range("D11") "0.07634"
Set-TextValue $ws.Range("E11") "  -1.88%  "
Set-TextValue $ws.Range("B12") "WrappedEther"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.764.35"
Set-TextValue $ws.Range("E12") "  +1.15%  "
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.570"
Set-TextValue $ws.Range("E13") "  -2.86%  "
Set-TextValue $ws.Range("D14") "1.948.72"
Set-TextValue $ws.Range("E14") "  -1.59%  "
Set-TextValue $ws.Range("D15") "0.5776"
Set-TextValue $ws.Range("E15") "  -3.51%  "
Set-TextValue $ws.Range("D16") "0.0₅8184"
Set-TextValue $ws.Range("E16") "  -2.80%  "
Set-TextValue $ws.Range("D17") "67.97"
Set-TextValue $ws.Range("E17") "  -1.62%  "
Set-TextValue $ws.Range("D18") "27.391.26"
Set-TextValue $ws.Range("E18") "  -1.57%  "
Set-TextValue $ws.Range("D19") "216.58"
Set-TextValue $ws.Range("E19") "  -4.00%  "
Set-TextValue $ws.Range("E20") "  +0.06%  "
Set-TextValue $ws.Range("E21") "  -3.51%  "
Set-TextValue $ws.Range("D22") "10.49"
Set-TextValue $ws.Range("E22") "  -4.01%  "
Set-TextValue $ws.Range("D23") "5.981"
Set-TextValue $ws.Range("E23") "  -4.27%  "
Set-TextValue $ws.Range("E24") "  +0.01%  "
Set-TextValue $ws.Range("D25") "142.49"
Set-TextValue $ws.Range("E25") "  -2.77%  "
Set-TextValue $ws.Range("D26") "1.735"
Set-TextValue $ws.Range("E26") "  +3.49%  "
Set-TextValue $ws.Range("D27") "0.1218"
Set-TextValue $ws.Range("E27") "  -2.67%  "
Set-TextValue $ws.Range("D28") "7.284"
Set-TextValue $ws.Range("E28") "  -2.51%  "
Set-TextValue $ws.Range("D29") "16.32"
Set-TextValue $ws.Range("E29") "  -5.18%  "
Set-TextValue $ws.Range("D30") "0.05429"
Set-TextValue $ws.Range("E30") "  -4.35%  "
Set-TextValue $ws.Range("E31") "  -1.40%  "
Set-TextValue $ws.Range("D32") "3.515"
Set-TextValue $ws.Range("E32") "  -5.00%  "
Set-TextValue $ws.Range("D33") "3.435"
Set-TextValue $ws.Range("E33") "  -2.66%  "
Set-TextValue $ws.Range("E34") "  -2.35%  "
Set-TextValue $ws.Range("D35") "2.880"
Set-TextValue $ws.Range("E35") "  +0.77%  "
Set-TextValue $ws.Range("D36") "0.9506"
Set-TextValue $ws.Range("E36") "  -2.82%  "
Set-TextValue $ws.Range("E37") "  -1.24%  "
Set-TextValue $ws.Range("D38") "0.5876"
Set-TextValue $ws.Range("E38") "  -1.74%  "
Set-TextValue $ws.Range("E39") "  -2.16%  "
Set-TextValue $ws.Range("D40") "5.879"
Set-TextValue $ws.Range("E40") "  -1.88%  "
Set-TextValue $ws.Range("D41") "1.044.18"
Set-TextValue $ws.Range("E41") "  -0.32%  "
Set-TextValue $ws.Range("E42") "  +0.02%  "
Set-TextValue $ws.Range("D43") "0.8413"
Set-TextValue $ws.Range("E43") "  -1.17%  "
Set-TextValue $ws.Range("D44") "101.01"
Set-TextValue $ws.Range("E44") "  -1.10%  "
Set-TextValue $ws.Range("D45") "1.854.55"
Set-TextValue $ws.Range("E45") "  -1.70%  "
Set-TextValue $ws.Range("E46") "  +8.91%  "
Set-TextValue $ws.Range("D47") "58.14"
Set-TextValue $ws.Range("E47") "  -2.39%  "
Set-TextValue $ws.Range("D48") "0.4512"
Set-TextValue $ws.Range("E48") "  +1.59%  "
Set-TextValue $ws.Range("E49") "  -0.24%  "
Set-TextValue $ws.Range("D50") "8.106"
Set-TextValue $ws.Range("E50") "  -2.38%  "
Set-TextValue $ws.Range("D51") "0.05241"
Set-TextValue $ws.Range("E51") "  -1.51%  "

Write-Host "Applied 93 cell updates"
